$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.2236597029243569; C = 0.6761979329963813; D = 0.6208552414790156; E = 0.7879436791287913; F = 0.7762373633456771; G = 19 }
    3  = @{ B = 0.1243738929906032; C = 0.5670298149555084; D = 0.4487500652448925; E = 0.6698880990470666; F = 0.677324390080757;  G = 18 }
    4  = @{ B = 0.1692601144216072; C = 0.4938470918437637; D = 0.3483999257463576; E = 0.5902541196352276; F = 0.5828682296416691; G = 17 }
    5  = @{ B = 0.2998499526236316; C = 0.5009693588263471; D = 0.3362117046188728; E = 0.5798376536746065; F = 0.5125639711936881; G = 16 }
    6  = @{ B = 0.3249303377857756; C = 0.510714471671429;  D = 0.3463237248587889; E = 0.5884927568447966; F = 0.5078779385610399; G = 15 }
    7  = @{ B = 0.357871698984652;  C = 0.5111964517621793; D = 0.3510459376640385; E = 0.592491297542874;  F = 0.4900261363687482; G = 14 }
    8  = @{ B = 0.3633620203708059; C = 0.5275006984044196; D = 0.3717216992278948; E = 0.6096898385473509; F = 0.5095722583647349; G = 13 }
    9  = @{ B = 0.4177687993451564; C = 0.536905318813658;  D = 0.3885450245204918; E = 0.623333798634802;  F = 0.4831874337779735; G = 12 }
    10 = @{ B = 0.435803605847665;  C = 0.5614457896225986; D = 0.409650469768552;  E = 0.6400394282921577; F = 0.4916281680178602; G = 11 }
    11 = @{ B = 0.4121559002167278; C = 0.5537197283926004; D = 0.4082939408543487; E = 0.638978826608792;  F = 0.5146967335462779; G = 10 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
